# "Allow blank banner. New changes to language captions"
#
# - Banner_Text!B4 (the M3 message) becomes the new combined Sanskrit
#   caption "श्री ब्रह्मतन्त्र  स्वतन्त्र परकाल स्वामि मठस्य आचार्याः"
#   (previously it held "श्रीपरकालमठस्य आचार्याः").
# - Banner_Text!B5 (the M4 message, previously
#   "ब्रह्मतन्त्रस्वतन्त्रस्वामिपरम्परा") is cleared out entirely so the
#   banner for that message id can be left blank.
# - Banner_Text!B6 (M5) keeps its existing text.
# - The active sheet moves from "acharyan_captions" to "Banner_Text", and
#   the last selections on "Founders_Early_Acharyas" / "Banner_Text" move
#   to D31 / B4 respectively.

$wb = $excel.ActiveWorkbook

$wsBanner   = $wb.Worksheets.Item("Banner_Text")
$wsFounders = $wb.Worksheets.Item("Founders_Early_Acharyas")

# New Sanskrit caption text for message M3.
$wsBanner.Range("B4").Value = "श्री ब्रह्मतन्त्र  स्वतन्त्र परकाल स्वामि मठस्य आचार्याः"

# Message M4 banner text is removed -> blank banner allowed.
$wsBanner.Range("B5").ClearContents()

# Leave the last selection on Founders_Early_Acharyas at D31.
$wsFounders.Range("D31").Select()

# Switch the active tab to Banner_Text and leave its selection on B4.
$wsBanner.Activate()
$wsBanner.Range("B4").Select()
